{"js": "// Lattice multiplication worksheet refresh: update the 5x3 table of exercises\n// in place. Each cell holds 5 lines (separated by manual line breaks):\n//   \"A x B\"\n//   \"  <d1(B)>    <d2(B)>\"\n//   \"  ----\"\n//   \"<d1(A)>|    |\"\n//   \"<d2(A)>|    |\"\n// The table shape (5 rows x 3 columns) does not change -- only the text.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst cellValues = [\n  { row: 0, col: 0, text: \"43 x 64\\v  6    4\\v  ----\\v4|    |\\v3|    |\" },\n  { row: 0, col: 1, text: \"73 x 28\\v  2    8\\v  ----\\v7|    |\\v3|    |\" },\n  { row: 0, col: 2, text: \"48 x 27\\v  2    7\\v  ----\\v4|    |\\v8|    |\" },\n  { row: 1, col: 0, text: \"73 x 81\\v  8    1\\v  ----\\v7|    |\\v3|    |\" },\n  { row: 1, col: 1, text: \"95 x 14\\v  1    4\\v  ----\\v9|    |\\v5|    |\" },\n  { row: 1, col: 2, text: \"48 x 79\\v  7    9\\v  ----\\v4|    |\\v8|    |\" },\n  { row: 2, col: 0, text: \"80 x 67\\v  6    7\\v  ----\\v8|    |\\v0|    |\" },\n  { row: 2, col: 1, text: \"84 x 74\\v  7    4\\v  ----\\v8|    |\\v4|    |\" },\n  { row: 2, col: 2, text: \"29 x 79\\v  7    9\\v  ----\\v2|    |\\v9|    |\" },\n  { row: 3, col: 0, text: \"20 x 74\\v  7    4\\v  ----\\v2|    |\\v0|    |\" },\n  { row: 3, col: 1, text: \"88 x 14\\v  1    4\\v  ----\\v8|    |\\v8|    |\" },\n  { row: 3, col: 2, text: \"82 x 15\\v  1    5\\v  ----\\v8|    |\\v2|    |\" },\n  { row: 4, col: 0, text: \"18 x 75\\v  7    5\\v  ----\\v1|    |\\v8|    |\" },\n  { row: 4, col: 1, text: \"71 x 83\\v  8    3\\v  ----\\v7|    |\\v1|    |\" },\n  { row: 4, col: 2, text: \"38 x 55\\v  5    5\\v  ----\\v3|    |\\v8|    |\" },\n];\n\nfor (const { row, col, text } of cellValues) {\n  const cell = table.getCell(row, col);\n  const range = cell.body.getRange(\"Whole\");\n  range.insertText(text, Word.InsertLocation.replace);\n}\n\nawait context.sync();", "ps1": "# Lattice multiplication worksheet refresh: update the 5x3 table of exercises\n# in place. Each cell holds 5 lines (separated by manual line breaks, char 11)\n# and the table shape (5 rows x 3 columns) is unchanged -- only the text.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$br = [char]11\n\n$cell = $t.Cell(1,1)\n$rng = $cell.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"43 x 64\" + $br + \"  6    4\" + $br + \"  ----\" + $br + \"4|    |\" + $br + \"3|    |\"\n\n$cell = $t.Cell(1,2)\n$rng = $cell.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"73 x 28\" + $br + \"  2    8\" + $br + \"  ----\" + $br + \"7|    |\" + $br + \"3|    |\"\n\n$cell = $t.Cell(1,3)\n$rng = $cell.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"48 x 27\" + $br + \"  2    7\" + $br + \"  ----\" + $br + \"4|    |\" + $br + \"8|    |\"\n\n$cell = $t.Cell(2,1)\n$rng = $cell.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"73 x 81\" + $br + \"  8    1\" + $br + \"  ----\" + $br + \"7|    |\" + $br + \"3|    |\"\n\n$cell = $t.Cell(2,2)\n$rng = $cell.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"95 x 14\" + $br + \"  1    4\" + $br + \"  ----\" + $br + \"9|    |\" + $br + \"5|    |\"\n\n$cell = $t.Cell(2,3)\n$rng = $cell.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"48 x 79\" + $br + \"  7    9\" + $br + \"  ----\" + $br + \"4|    |\" + $br + \"8|    |\"\n\n$cell = $t.Cell(3,1)\n$rng = $cell.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"80 x 67\" + $br + \"  6    7\" + $br + \"  ----\" + $br + \"8|    |\" + $br + \"0|    |\"\n\n$cell = $t.Cell(3,2)\n$rng = $cell.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"84 x 74\" + $br + \"  7    4\" + $br + \"  ----\" + $br + \"8|    |\" + $br + \"4|    |\"\n\n$cell = $t.Cell(3,3)\n$rng = $cell.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"29 x 79\" + $br + \"  7    9\" + $br + \"  ----\" + $br + \"2|    |\" + $br + \"9|    |\"\n\n$cell = $t.Cell(4,1)\n$rng = $cell.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"20 x 74\" + $br + \"  7    4\" + $br + \"  ----\" + $br + \"2|    |\" + $br + \"0|    |\"\n\n$cell = $t.Cell(4,2)\n$rng = $cell.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"88 x 14\" + $br + \"  1    4\" + $br + \"  ----\" + $br + \"8|    |\" + $br + \"8|    |\"\n\n$cell = $t.Cell(4,3)\n$rng = $cell.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"82 x 15\" + $br + \"  1    5\" + $br + \"  ----\" + $br + \"8|    |\" + $br + \"2|    |\"\n\n$cell = $t.Cell(5,1)\n$rng = $cell.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"18 x 75\" + $br + \"  7    5\" + $br + \"  ----\" + $br + \"1|    |\" + $br + \"8|    |\"\n\n$cell = $t.Cell(5,2)\n$rng = $cell.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"71 x 83\" + $br + \"  8    3\" + $br + \"  ----\" + $br + \"7|    |\" + $br + \"1|    |\"\n\n$cell = $t.Cell(5,3)\n$rng = $cell.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"38 x 55\" + $br + \"  5    5\" + $br + \"  ----\" + $br + \"3|    |\" + $br + \"8|    |\"\n"}
